# "Actualizar" refresh: the availability-check timestamps in column D
# (Fecha) roll down one block and the newest block gets a fresh check
# timestamp, matching the automated "Actualizar 02-13-2021 13-15-02" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest batch (rows 2-15) recorded a brand-new timestamp.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 44240.55202552549
}

# Previous newest batch (rows 16-29) shifts down, taking on what used
# to be the rows 2-15 timestamp (re-measured on write, hence the tiny
# float delta vs. the old value).
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 44240.53084094908
}

# Oldest batch (rows 30-43) shifts down, taking the previous rows
# 16-29 timestamp verbatim.
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 44240.5096403125
}
